$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "26_03_2024"
$ws.Range("F2").Value = 4247
$ws.Range("F3").Value = 3546
$ws.Range("F4").Value = 4411
$ws.Range("F5").Value = 548

$ws.Range("F6").Select()
